# "Alteração ppt com componentes"
#
# 1) The automatic date/time footer field (type="datetimeFigureOut") that
#    appears on the Slide Master and on every Slide Layout is bumped from
#    27/11/2021 -> 29/11/2021.
# 2) The Slide Master background switches from the theme's bg1 scheme
#    color reference to an explicit solid fill (RGB 1B2024).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

$oldDate = "27/11/2021"
$newDate = "29/11/2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try {
            $phType = $shp.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Update the date placeholder on the Slide Master itself.
Update-DatePlaceholder $master.Shapes

# Update the date placeholder on every Slide Layout belonging to the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide Master background: explicit solid fill instead of the bg1 theme ref.
# RGB(0x1B, 0x20, 0x24) packed as 0x00BBGGRR (OLE COLORREF byte order).
$master.Background.Fill.ForeColor.RGB = 0x24201B
